$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.589.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.466.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.462.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.577"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.437"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.063.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.628.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.472.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.530"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -5.65%  "
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.818"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.850.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("E44").Value = "  -4.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("E47").Value = "  +11.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "333.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.47%  "
$ws.Range("E49").Value = "  -2.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.845"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.05%  "
